$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 270, shifting existing rows 270-319 down to 271-320
$ws.Rows.Item(270).Insert()

# Populate the new row 270 with the new observation
$ws.Cells.Item(270, 1).Value = 8
$ws.Cells.Item(270, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(270, 3).Value = "Coquimbo"
$ws.Cells.Item(270, 4).Value = 44694
$ws.Cells.Item(270, 5).Value = 4
$ws.Cells.Item(270, 6).Value = 100112032
$ws.Cells.Item(270, 7).Value = "Zapallo italiano"
$ws.Cells.Item(270, 8).Value = "Sin especificar"
$ws.Cells.Item(270, 9).Value = "Primera"
$ws.Cells.Item(270, 10).Value = 440
$ws.Cells.Item(270, 11).Value = 15000
$ws.Cells.Item(270, 12).Value = 16000
$ws.Cells.Item(270, 13).Value = 15500
$ws.Cells.Item(270, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(270, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(270, 16).Value = 258
$ws.Cells.Item(270, 17).Value = 60
$ws.Cells.Item(270, 18).Value = "Hortaliza"
